$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I3").Value = "Start of week 6 is Feb 17"
